# "adding create new product"
# 1) Update the client rep's and company rep's e-mail addresses (domain changed
#    from fai.ws -> fai.sa) on the "clients" and "company" sheets.
# 2) Add a brand-new "products" sheet at the end of the workbook with a
#    header row + one sample/auto-generated product row, and make it the
#    active sheet.

$wb = $excel.ActiveWorkbook

# --- clients sheet: repEmail (X2) and mail (S2) ---------------------------
$clients = $wb.Worksheets.Item("clients")
$clients.Range("X2").Value = "mario@fai.sa"
$clients.Range("S2").Value = "menna@fai.sa"
$clients.Columns.Item(24).ColumnWidth = 17.3854167
$null = $clients.Range("X2").Select()

# --- company sheet: repEmail (J2) ------------------------------------------
$company = $wb.Worksheets.Item("company")
$company.Range("J2").Value = "mario@fai.sa"
$null = $company.Range("J2").Select()

# --- purchaseCashback sheet: select whole column F (no longer the tab shown)
$purchaseCashback = $wb.Worksheets.Item("purchaseCashback")
$null = $purchaseCashback.Range("F1:F1048576").Select()

# --- new "products" sheet ---------------------------------------------------
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$products = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $last)
$products.Name = "products"

$products.Range("A1").Value = "ArName"
$products.Range("A2").Value = "منتج"
$products.Range("B1").Value = "EnName"
$products.Range("B2").Value = "product"
$products.Range("C1").Value = "ArDescription"
$products.Range("C2").Value = "وصف المنتج"
$products.Range("D1").Value = "EnDescription"
$products.Range("D2").Value = "product description"
$products.Range("E1").Value = "productType"
$products.Range("E2").Value = "كلاهما"
$products.Range("F1").Value = "productPrice"
$products.Range("F2").Value = 1000
$products.Range("G1").Value = "productCode"
$products.Range("G2").Value = "Auto"

$products.Columns.Item(1).ColumnWidth = 11.49869795
$products.Columns.Item(2).ColumnWidth = 10.7213542
$products.Columns.Item(3).ColumnWidth = 16.1666667
$products.Columns.Item(4).ColumnWidth = 18.1666667
$products.Columns.Item(5).ColumnWidth = 12.2760417
$products.Columns.Item(6).ColumnWidth = 13.8307292
$products.Columns.Item(7).ColumnWidth = 10.7213542

$null = $products.Range("G7").Select()
